$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.284.13'
$ws.Range("E2").Value = '  +0.01%  '

$ws.Range("D3").Value = '1.931.56'
$ws.Range("E3").Value = '  +0.05%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7516'
$ws.Range("E5").Value = '  +4.87%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.94'
$ws.Range("E6").Value = '  -2.80%  '

$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '27.76'
$ws.Range("E8").Value = '  +0.20%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3172'
$ws.Range("E9").Value = '  -0.64%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07094'
$ws.Range("E10").Value = '  -0.20%  '

$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08037'
$ws.Range("E11").Value = '  +0.74%  '

$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7778'
$ws.Range("E12").Value = '  -1.68%  '

$ws.Range("D13").Value = '1.920.83'
$ws.Range("E13").Value = '  -0.43%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.383'
$ws.Range("E14").Value = '  -0.06%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.04'
$ws.Range("E15").Value = '  -1.88%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.54'
$ws.Range("E16").Value = '  -0.96%  '

$ws.Range("D17").Value = '30.277.07'
$ws.Range("E17").Value = '  +0.01%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.991'
$ws.Range("E18").Value = '  +3.81%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '251.56'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007938'
$ws.Range("E20").Value = '  -1.30%  '

$ws.Range("D21").Value = '2.179.92'
$ws.Range("E21").Value = '  +0.16%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9998'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.05%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.680'
$ws.Range("E24").Value = '  -2.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.519'
$ws.Range("E25").Value = '  -0.23%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.75'
$ws.Range("E26").Value = '  -0.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.08'
$ws.Range("E27").Value = '  -0.16%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1298'
$ws.Range("E28").Value = '  +2.34%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.186'
$ws.Range("E29").Value = '  -3.45%  '

$ws.Range("E30").Value = '  +0.91%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.543'
$ws.Range("E31").Value = '  +1.06%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.407'
$ws.Range("E32").Value = '  +0.29%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.134'
$ws.Range("E33").Value = '  -0.01%  '

$ws.Range("E34").Value = '  +1.60%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.315'
$ws.Range("E35").Value = '  +3.42%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7565'
$ws.Range("E36").Value = '  +1.56%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.779'
$ws.Range("E37").Value = '  +0.40%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01951'
$ws.Range("E38").Value = '  -0.43%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.797'
$ws.Range("E39").Value = '  -0.07%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '78.22'
$ws.Range("E40").Value = '  -0.37%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.496'
$ws.Range("E41").Value = '  +2.17%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4517'
$ws.Range("E42").Value = '  +0.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.976'
$ws.Range("E43").Value = '  -0.60%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8404'
$ws.Range("E44").Value = '  -0.72%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9999'
$ws.Range("E45").Value = '  +0.06%  '

$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.660'
$ws.Range("E46").Value = '  +3.30%  '

$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.65'
$ws.Range("E47").Value = '  +1.13%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.936'
$ws.Range("E48").Value = '  +1.60%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '38.02'
$ws.Range("E49").Value = '  +3.76%  '

$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '961.24'
$ws.Range("E50").Value = '  +1.54%  '

$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1216'
$ws.Range("E51").Value = '  +6.52%  '
